$wb = $excel.ActiveWorkbook

# --- Reorder worksheet tabs: review_info should come before hotel_info ---
$reviewSheet = $wb.Worksheets.Item("review_info")
$hotelSheet = $wb.Worksheets.Item("hotel_info")
$reviewSheet.Move($hotelSheet)

# Sheet object references in this runtime are bound by tab index, so after the
# Move() above the old $hotelSheet variable no longer points at "hotel_info".
# Re-resolve it by name before making further edits.
$hotelSheet = $wb.Worksheets.Item("hotel_info")

# --- Add a new "State" column to hotel_info (inserted right after Hotel_Name, before City) ---
$hotelSheet.Columns.Item(3).Insert()
$hotelSheet.Cells.Item(1, 3).Value = "State"
$hotelSheet.Cells.Item(2, 3).Value = "Louisiana"
